# "added files for 4/22" - log hours worked on the "TA Meeting" and
# "Audio Backend" task rows of the timesheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# TA Meeting row (row 10): 1 hour logged
$ws.Range("E10").Value = 1

# Audio Backend row (row 13): hours logged across the week
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 2
$ws.Range("H13").Value = 1

# Leave the sheet scrolled/selected where the author last left off
$ws.Activate()
$ws.Range("D9").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
